# Actualizacion desde MV -datos-
# Appends 5 new daily rows (04-10-2021 .. 08-10-2021) to the "Diaria" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 194

$data = @(
    @("04-10-2021", 4.84, 4.08, 3.51, 3.37),
    @("05-10-2021", 4.83, 4.09, 3.53, 3.36),
    @("06-10-2021", 4.91, 4.1,  3.54, 3.37),
    @("07-10-2021", 5.03, 4.16, 3.57, 3.37),
    @("08-10-2021", 5.51, 4.49, 3.73, 3.42)
)

$r = $startRow
foreach ($row in $data) {
    $dateText = $row[0]

    # Write the date as a text formula literal and flatten it to a value.
    # This keeps the cell as a shared-string ("t=s") entry, matching the
    # existing column-A cells, instead of letting Excel auto-convert the
    # "dd-mm-yyyy" looking text into a date serial number (which would
    # also introduce a brand-new number format / cell style).
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = "=""" + $dateText + """"
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)  # xlPasteValues

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    $r = $r + 1
}

$excel.CutCopyMode = 0
